$d = $word.ActiveDocument

$rng = $d.Content.Duplicate
$found = $rng.Find.Execute("Unraveling the Enigma of Consciousness", $true, $false, $false, $false, $false, $true, 1, $false)
if (-not $found) { throw "Not found: Unraveling the Enigma of Consciousness" }
$rng.Text = "Arts and Expression: Exploring the Human Palette"

$rng = $d.Content.Duplicate
$found = $rng.Find.Execute("Emily Carter", $true, $false, $false, $false, $false, $true, 1, $false)
if (-not $found) { throw "Not found: Emily Carter" }
$rng.Text = "Sophia Mitchell"

$rng = $d.Content.Duplicate
$found = $rng.Find.Execute("emily", $true, $false, $false, $false, $false, $true, 1, $false)
if (-not $found) { throw "Not found: emily" }
$rng.Text = "sophiamitchell"

$rng = $d.Content.Duplicate
$found = $rng.Find.Execute("carter@cognit-research", $true, $false, $false, $false, $false, $true, 1, $false)
if (-not $found) { throw "Not found: carter@cognit-research" }
$rng.Text = "arts@educonnect"

$rng = $d.Content.Duplicate
$found = $rng.Find.Execute("From the dawn of human civilization, consciousness has stood as an enigmatic puzzle, captivating the minds of philosophers, scientists, and artists alike", $true, $false, $false, $false, $false, $true, 1, $false)
if (-not $found) { throw "Not found: From the dawn of human civilization, consciousness has stood as an enigmatic puzzle, captivating the minds of philosophers, scientists, and artists alike" }
$rng.Text = "In the realm of human creativity, arts serve as a captivating mirror to our souls, allowing us to delve into the labyrinth of emotions, ideas, and experiences that shape our journey through life"

$rng = $d.Content.Duplicate
$found = $rng.Find.Execute(" What is the nature of consciousness? How does it arise from the intricate workings of the brain? As we embark on this exploration, we will traverse diverse disciplines, weaving together insights from neuroscience, psychology, philosophy, and even quantum physics, in a quest to unravel the enigma of consciousness", $true, $false, $false, $false, $false, $true, 1, $false)
if (-not $found) { throw "Not found:  What is the nature of consciousness? How does it arise from the intricate workings of the brain? As we embark on this exploration, we will traverse diverse disciplines, weaving together insights from neuroscience, psychology, philosophy, and even quantum physics, in a quest to unravel the enigma of consciousness" }
$rng.Text = " Like a symphony of colors on an artist's canvas, arts resonate with the uniqueness of every individual, unveiling the kaleidoscope of perspectives that define our collective human tapestry. Through the prism of diverse artistic mediums, be it the written word, visual artistry, or the emotive power of music, we embark on a quest to explore the profound impact arts have on our perception of the world around us"

$rng = $d.Content.Duplicate
$found = $rng.Find.Execute("Seeking answers, we delve into the depths of neuroscience, charting the intricate neural networks that serve as the physical substrate of consciousness", $true, $false, $false, $false, $false, $true, 1, $false)
if (-not $found) { throw "Not found: Seeking answers, we delve into the depths of neuroscience, charting the intricate neural networks that serve as the physical substrate of consciousness" }
$rng.Text = "Arts possess a formidable power to transcend linguistic boundaries, breaking down the walls that separate cultures and bringing hearts closer"

$rng = $d.Content.Duplicate
$found = $rng.Find.Execute(" We scrutinize the interplay of neurons, synapses, and brain regions, seeking to decipher how these biological components orchestrate the symphony of conscious experience", $true, $false, $false, $false, $false, $true, 1, $false)
if (-not $found) { throw "Not found:  We scrutinize the interplay of neurons, synapses, and brain regions, seeking to decipher how these biological components orchestrate the symphony of conscious experience" }
$rng.Text = " Like a dancer's graceful movements, arts have the ability to bypass words and communicate emotions in a universal language understood by all"

$rng = $d.Content.Duplicate
$found = $rng.Find.Execute(" We contemplate the role of attention, memory, and emotion, probing their contribution to the subjective tapestry of consciousness", $true, $false, $false, $false, $false, $true, 1, $false)
if (-not $found) { throw "Not found:  We contemplate the role of attention, memory, and emotion, probing their contribution to the subjective tapestry of consciousness" }
$rng.Text = " Whether it's the haunting strains of a heartfelt melody, the evocative imagery of a painting, or the poignant verses of a poem, arts have a remarkable capacity to unify humanity, fostering a sense of connectedness and empathy that resounds across geographical and cultural divides"

$rng = $d.Content.Duplicate
$found = $rng.Find.Execute("Venturing beyond the confines of neuroscience, we turn to psychology, seeking insights into the phenomenology of consciousness", $true, $false, $false, $false, $false, $true, 1, $false)
if (-not $found) { throw "Not found: Venturing beyond the confines of neuroscience, we turn to psychology, seeking insights into the phenomenology of consciousness" }
$rng.Text = "Beyond their aesthetic appeal, arts have a transformative influence on the human psyche, playing a pivotal role in personal development and overall well-being"

$rng = $d.Content.Duplicate
$found = $rng.Find.Execute(" Through introspection and experimentation, we dissect the various dimensions of conscious experience, from the vivid hues of color to the ethereal flow of time", $true, $false, $false, $false, $false, $true, 1, $false)
if (-not $found) { throw "Not found:  Through introspection and experimentation, we dissect the various dimensions of conscious experience, from the vivid hues of color to the ethereal flow of time" }
$rng.Text = " Engaging with arts in any form, whether as a creator or an audience, provides an avenue for self-expression, emotional release, and introspection"

$rng = $d.Content.Duplicate
$found = $rng.Find.Execute(" We explore altered states of consciousness, such as dreams, meditation, and psychedelic experiences, seeking clues to the malleability and multidimensionality of our inner worlds", $true, $false, $false, $false, $false, $true, 1, $false)
if (-not $found) { throw "Not found:  We explore altered states of consciousness, such as dreams, meditation, and psychedelic experiences, seeking clues to the malleability and multidimensionality of our inner worlds" }
$rng.Text = " Just as a sculptor chisels away at a block of marble, revealing the hidden beauty within, arts empower us to confront our own complexities, unraveling the enigmas of our inner selves"

$rng = $d.Content.Duplicate
$found = $rng.Find.Execute("In this exploration of consciousness, we embarked on an interdisciplinary voyage, delving into the depths of neuroscience, psychology, philosophy, and quantum physics", $true, $false, $false, $false, $false, $true, 1, $false)
if (-not $found) { throw "Not found: In this exploration of consciousness, we embarked on an interdisciplinary voyage, delving into the depths of neuroscience, psychology, philosophy, and quantum physics" }
$rng.Text = "In the realm of arts, we find a mirror to humanity's soul, a kaleidoscope of perspectives, and a profound force for connection and transformation"

$rng = $d.Content.Duplicate
$found = $rng.Find.Execute(" Neuroscience provided invaluable insights into the neural underpinnings of consciousness, revealing the intricate dance of neurons, synapses, and brain regions that orchestrate our subjective experiences", $true, $false, $false, $false, $false, $true, 1, $false)
if (-not $found) { throw "Not found:  Neuroscience provided invaluable insights into the neural underpinnings of consciousness, revealing the intricate dance of neurons, synapses, and brain regions that orchestrate our subjective experiences" }
$rng.Text = " Through the diverse mediums of expression, arts transcend boundaries, fostering empathy and unifying people from all walks of life"

$rng = $d.Content.Duplicate
$found = $rng.Find.Execute(" Psychology illuminated the phenomenology of consciousness, shedding light on the kaleidoscope of subjective experiences that comprise our inner worlds. Philosophy and quantum physics challenged our conventional notions of consciousness, inviting us to ponder the possibility of alternative frameworks and unexplored dimensions. As our understanding of ", $true, $false, $false, $false, $false, $true, 1, $false)
if (-not $found) { throw "Not found:  Psychology illuminated the phenomenology of consciousness, shedding light on the kaleidoscope of subjective experiences that comprise our inner worlds. Philosophy and quantum physics challenged our conventional notions of consciousness, inviting us to ponder the possibility of alternative frameworks and unexplored dimensions. As our understanding of " }
$rng.Text = " From the depths of introspection to the heights of collective celebration, arts serve as a catalyst for "

$rng = $d.Content.Duplicate
$found = $rng.Find.Execute("consciousness continues to evolve, we remain humbled by its enigmatic nature, yet driven by an insatiable curiosity to unravel its mysteries", $true, $false, $false, $false, $false, $true, 1, $false)
if (-not $found) { throw "Not found: consciousness continues to evolve, we remain humbled by its enigmatic nature, yet driven by an insatiable curiosity to unravel its mysteries" }
$rng.Text = "self-discovery, emotional release, and the cultivation of a more profound understanding of ourselves and the world around us"

